$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated financial figures per row (columns D..AJ), keyed by cell address.
$updates = @{
    "D2" = 3841
    "E2" = 178
    "F2" = 178
    "G2" = 159
    "H2" = 0
    "I2" = 24
    "J2" = -24
    "K2" = 4375
    "L2" = 2434
    "M2" = 1941
    "N2" = 1453
    "O2" = 488
    "P2" = 112
    "Q2" = 268
    "R2" = -338
    "S2" = -10
    "T2" = 334
    "U2" = -66
    "V2" = 850
    "W2" = 4.63
    "X2" = 0.01
    "Y2" = 1.67
    "Z2" = 0.01
    "AA2" = 125.41
    "AB2" = 1213.86
    "AC2" = 109
    "AD2" = 96.7
    "AE2" = 6569
    "AF2" = 1.61
    "AG2" = 50
    "AH2" = 0.47
    "AI2" = 45.2
    "AJ2" = 22434980
    "D3" = 3777
    "E3" = 161
    "F3" = 161
    "G3" = 126
    "H3" = 134
    "I3" = 136
    "J3" = -2
    "K3" = 4503
    "L3" = 2545
    "M3" = 1958
    "N3" = 1539
    "O3" = 419
    "P3" = 112
    "Q3" = 390
    "R3" = -393
    "S3" = -40
    "T3" = 338
    "U3" = 52
    "V3" = 892
    "W3" = 4.27
    "X3" = 3.55
    "Y3" = 9.09
    "Z3" = 3.02
    "AA3" = 129.96
    "AB3" = 1285.44
    "AC3" = 606
    "AD3" = 15.1
    "AE3" = 6955
    "AF3" = 1.32
    "AG3" = 80
    "AH3" = 0.87
    "AI3" = 13.01
    "AJ3" = 22434980
    "D4" = 3914
    "E4" = 259
    "F4" = 259
    "G4" = 172
    "H4" = 67
    "I4" = 47
    "J4" = 21
    "K4" = 4628
    "L4" = 2544
    "M4" = 2083
    "N4" = 1580
    "O4" = 503
    "P4" = 112
    "Q4" = 350
    "R4" = -410
    "S4" = 9
    "T4" = 173
    "U4" = 178
    "V4" = 912
    "W4" = 6.61
    "X4" = 1.72
    "Y4" = 2.99
    "Z4" = 1.48
    "AA4" = 122.13
    "AB4" = 1310.69
    "AC4" = 208
    "AD4" = 36.23
    "AE4" = 7143
    "AF4" = 1.06
    "AG4" = 100
    "AH4" = 1.33
    "AI4" = 47.38
    "AJ4" = 22434980
    "D5" = 3524
    "E5" = 198
    "F5" = 198
    "G5" = 176
    "H5" = 94
    "I5" = 51
    "J5" = 42
    "K5" = 5199
    "L5" = 2718
    "M5" = 2482
    "N5" = 1621
    "O5" = 860
    "P5" = 112
    "Q5" = 146
    "R5" = 274
    "S5" = -171
    "T5" = 113
    "U5" = 33
    "V5" = 788
    "W5" = 5.6
    "X5" = 2.66
    "Y5" = 3.21
    "Z5" = 1.91
    "AA5" = 109.51
    "AB5" = 1347.76
    "AC5" = 229
    "AD5" = 22.73
    "AE5" = 7327
    "AF5" = 0.71
    "AG5" = 110
    "AH5" = 2.12
    "AI5" = 47.42
    "AJ5" = 22434980
    "D6" = 5347
    "E6" = 133
    "F6" = 133
    "G6" = 224
    "H6" = 188
    "I6" = 128
    "K6" = 6708
    "L6" = 3841
    "M6" = 2867
    "N6" = 1748
    "P6" = 112
    "Q6" = 14
    "R6" = -114
    "S6" = 216
    "T6" = 224
    "U6" = -211
    "V6" = 1126
    "W6" = 2.48
    "X6" = 3.51
    "Y6" = 7.6
    "Z6" = 3.15
    "AA6" = 133.96
    "AB6" = 1429.94
    "AC6" = 571
    "AD6" = 6.52
    "AE6" = 7902
    "AF6" = 0.47
    "AG6" = 80
    "AH6" = 2.15
    "AI6" = 13.82
    "AJ6" = 22434980
    "D7" = 6565
    "E7" = 187
    "I7" = 121
    "W7" = 2.85
    "AC7" = 539
    "AD7" = 9.01
    "D8" = 7006
    "E8" = 253
    "I8" = 120
    "W8" = 3.61
    "AC8" = 535
    "AD8" = 9.09
    "D9" = 7539
    "E9" = 304
    "I9" = 147
    "W9" = 4.03
    "AC9" = 655
    "AD9" = 7.42
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# Cells that no longer have reported figures for rows 7-9 (2019E-2021E) are cleared.
$clears = @(
    "G7", "H7", "K7", "L7", "M7", "N7", "P7", "Q7",
    "R7", "S7", "T7", "U7", "X7", "Y7", "Z7", "AA7",
    "AE7", "AF7", "AG7", "AH7", "AI7", "G8", "H8", "K8",
    "L8", "M8", "N8", "P8", "Q8", "R8", "S8", "T8",
    "U8", "X8", "Y8", "Z8", "AA8", "AE8", "AF8", "AG8",
    "AH8", "AI8", "G9", "H9", "K9", "L9", "M9", "N9",
    "P9", "Q9", "R9", "S9", "T9", "U9", "X9", "Y9",
    "Z9", "AA9", "AE9", "AF9", "AG9", "AH9", "AI9"
)

foreach ($addr in $clears) {
    $ws.Range($addr).ClearContents()
}
